# Update "想去人数" (wish-to-go count) figures in column F across sheets
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 309  # was 307
$ws1.Range("F4").Value = 186  # was 185
$ws1.Range("F5").Value = 1182  # was 1178
$ws1.Range("F8").Value = 128  # was 126
$ws1.Range("F9").Value = 135  # was 132
$ws1.Range("F11").Value = 256  # was 255
$ws1.Range("F12").Value = 153  # was 150
$ws1.Range("F13").Value = 155  # was 151
$ws1.Range("F14").Value = 1372  # was 1363
$ws1.Range("F15").Value = 522  # was 517
$ws1.Range("F16").Value = 198  # was 196
$ws1.Range("F17").Value = 311  # was 308
$ws1.Range("F19").Value = 716  # was 711
$ws1.Range("F20").Value = 1099  # was 1098
$ws1.Range("F22").Value = 1926  # was 1928
$ws1.Range("F23").Value = 2551  # was 2547
$ws1.Range("F24").Value = 1311  # was 1306
$ws1.Range("F26").Value = 222  # was 219
$ws1.Range("F27").Value = 376  # was 375
$ws1.Range("F28").Value = 955  # was 944
$ws1.Range("F30").Value = 1101  # was 1091
$ws1.Range("F31").Value = 133  # was 132
$ws1.Range("F33").Value = 758  # was 757
$ws1.Range("F34").Value = 445  # was 437
$ws1.Range("F35").Value = 598  # was 596
$ws1.Range("F36").Value = 767  # was 762
$ws1.Range("F37").Value = 327  # was 326
$ws1.Range("F38").Value = 217  # was 216

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 176  # was 174
$ws2.Range("F13").Value = 532  # was 526

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 309  # was 307
$ws4.Range("F7").Value = 186  # was 185
$ws4.Range("F9").Value = 176  # was 174
$ws4.Range("F10").Value = 1182  # was 1178
$ws4.Range("F13").Value = 128  # was 126
$ws4.Range("F15").Value = 135  # was 132
$ws4.Range("F17").Value = 256  # was 255
$ws4.Range("F19").Value = 153  # was 150
$ws4.Range("F20").Value = 155  # was 151
$ws4.Range("F21").Value = 1372  # was 1363
$ws4.Range("F22").Value = 522  # was 517
$ws4.Range("F23").Value = 198  # was 196
$ws4.Range("F24").Value = 311  # was 308
$ws4.Range("F26").Value = 1099  # was 1098
$ws4.Range("F27").Value = 2551  # was 2547
$ws4.Range("F29").Value = 1311  # was 1306
$ws4.Range("F34").Value = 222  # was 219
$ws4.Range("F35").Value = 376  # was 375
$ws4.Range("F36").Value = 955  # was 944
$ws4.Range("F40").Value = 1101  # was 1091
$ws4.Range("F41").Value = 758  # was 757
$ws4.Range("F42").Value = 445  # was 437
$ws4.Range("F43").Value = 598  # was 596
$ws4.Range("F44").Value = 767  # was 762
$ws4.Range("F45").Value = 327  # was 326
$ws4.Range("F48").Value = 217  # was 216

Write-Output "Done updating F column values."
